$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Move the merged title ("Data Siswa ...") block up one row: it used to
#    live in A2 (merged A2:Q3); it should now live in A1 (merged A1:Q3) and
#    every cell in the A1:Q3 block should carry the centered/word-wrapped
#    style that the title cell already used.
# ---------------------------------------------------------------------------

$titleText = $ws.Range("A2").Value()

# Unmerge the old title block so we can freely move values/styles around.
$ws.Range("A2:Q3").UnMerge()

# Copy the title cell's format (centered + wrap text) onto the whole
# A1:Q3 block first ...
$ws.Range("A2").Copy()
$ws.Range("A1:Q3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ... then clear out the old text in A2 and write it into A1 instead.
$ws.Range("A2").ClearContents()
$ws.Range("A1").Value = $titleText

# Re-merge the block, now anchored at A1.
$ws.Range("A1:Q3").Merge()

# ---------------------------------------------------------------------------
# 2. Sheet view: scroll the visible window so column K is the left-most
#    visible column (was column G).
# ---------------------------------------------------------------------------

$win = $wb.Application.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1

Write-Host "done"
